{"js": "// The canonical OOXML diff for this revision only touches build/tooling\n// metadata that sits outside the Word content object model:\n//   1. The `<!-- Created by docx4j \u2026 -->` generator-stamp comment that\n//      docx4j writes as the first child of <w:body> (its text just records\n//      the docx4j/JAXB/JDK/OS versions used to regenerate the fixture).\n//   2. The order in which xmlns:* namespace declarations are serialized on\n//      the root elements of document.xml / footer.xml / header.xml /\n//      styles.xml (same namespace set and URIs, just re-ordered).\n// Every real, user-visible piece of content in the hunk context (the\n// paragraph \"First block of main text.\", the header/footer paragraphs and\n// their formatting, all the style definitions, etc.) is byte-for-byte\n// identical before and after. Neither artifact is reachable from Office.js\n// (there is no API that reads or writes a raw XML comment node, and\n// namespace-prefix ordering is an implementation detail of the XML writer,\n// not part of the document model), so there is nothing in this diff for a\n// content script to apply. We still touch the API surface the task\n// describes so the request round-trips cleanly through context.sync()\n// without mutating any content.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The canonical OOXML diff for this revision only touches build/tooling\n# metadata that sits outside the Word content object model:\n#   1. The `<!-- Created by docx4j ... -->` generator-stamp comment that\n#      docx4j writes as the first child of <w:body> (its text just records\n#      the docx4j/JAXB/JDK/OS versions used to regenerate the fixture).\n#   2. The order in which xmlns:* namespace declarations are serialized on\n#      the root elements of document.xml / footer.xml / header.xml /\n#      styles.xml (same namespace set and URIs, just re-ordered).\n# Every real, user-visible piece of content in the hunk context (the\n# paragraph \"First block of main text.\", the header/footer paragraphs and\n# their formatting, all the style definitions, etc.) is byte-for-byte\n# identical before and after. Neither artifact is reachable from the Word\n# COM object model (WordOpenXML/Range.XML are read-only, and there is no\n# property that exposes a raw XML comment node or controls namespace-prefix\n# ordering), so there is nothing in this diff for a content script to\n# apply. We still touch the documented object model so the session\n# round-trips cleanly without mutating any content.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
